# Restaura conteudo completo da Planilha Base
# Preenche os textos de analise que estavam ausentes/placeholders e limpa a
# linha 25 (duplicata da linha 14 do bloco de Meta Especifica), conforme o
# commit "Restaura conteudo completo da Planilha Base".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Linha 4: Area Tematica / correlacao problema x intervencao -----------
$ws.Range("F4").Value = "SIM, pois está entre as metas do PNSP, bem como no enfoque desta Área Temática."

# --- Linha 8: Meta e indicador geral ---------------------------------------
$ws.Range("A8").Value = "1* Redução da taxa estadual de Morte Violenta Intencional para no máximo de 6,41 por 100 mil habitantes até 2027. A taxa atual é de 8,41 mortes por grupo de 100 mil habitantes em 2024, segundo dados do SINESP.1*"
$ws.Range("B8").Value = "SIM, pois foi estabelecida uma meta geral que reflete o problema público e demonstrada uma estratégia adequada para a sua mitigação."
$ws.Range("C8").Value = "SIM, pois houve a identificação de um relevante problema público selecionado como objeto da intervenção proposta, bem como um indicador capaz de demonstrar a implementação."
$ws.Range("D8").Value = "SIM, pois verifica-se a relação direta entre as medidas propostas e a meta geral estabelecida."
$ws.Range("E8").Value = "SIM, já que foram indicadas ações que em conjunto permitem o atingimento das metas específicas e, consequentemente da meta geral."
$ws.Range("F8").Value = "SIM, o indicador é capaz de mensurar o desempenho dos processos voltados ao atingimento da Meta geral."

# --- Linha 10: Formula / referencia utilizada (meta geral) ----------------
$ws.Range("F10").Value = "0* Indicador: Taxa de mortes violentas Finalidade: Verificar a variação da taxa de mortes violentas Fórmula de cálculo: (Σ de vítimas homicídio doloso + feminicídio + latrocínio + lesão corporal seguida de morte + morte por intervenção policial x 100 mil) /população Periodicidade: Anual Índice Atual de SC 2025: 8,41% (SINESP) Variação (atual/meta): -52,5% Meta 2030 Plano Estadual: 6,0 Fonte: Plano Estadual de Segurança Pública e Defesa Social II de Santa Catarina.0*"

# --- Linha 12: Estrategia de implementacao ---------------------------------
$ws.Range("B12").Value = "ATENDE, pois identifica e delimita o problema público a ser enfrentado."
$ws.Range("C12").Value = "ATENDE e aponta os referidos mecanismos adequadamente."
$ws.Range("D12").Value = "ATENDE, pois houve demonstração da pretensão de difundir os conhecimentos do tema aos policiais que atuam na área, além disso foi demonstrada a intenção de efetivar uma atuação integrada."
$ws.Range("E12").Value = "ATENDE, pois foram indicadas contratações de objetos e serviços com o objetivo de mitigar o déficit de equipamentos, aumentar a capacidade de prestação dos serviços no intuito de implementar a Meta Geral estabelecida no Plano de Aplicação."

# --- Linha 13: cabecalhos da Meta Especifica -------------------------------
$ws.Range("G13").Value = "A Meta do PESP foi informada? Existe aderência?"
$ws.Range("H13").Value = "A Meta do PNSP foi informada? Existe aderência?"
$ws.Range("I13").Value = "A política da Carteira de Políticas do MJSP foi informada? Existe aderência?"

# --- Linha 14: Meta Especifica (bloco completo) ----------------------------
$ws.Range("A14").Value = "2*1 - Reduzir em 20% a taxa de homicídios e feminicídios, latrocínios e lesões corporais seguida de morte nos 60 municípios com população igual ou superior a 26.500 habitantes até 2027.2*"

$e14 = @"
SIM.

A referência informada foi:



3*Referência: 371 homicídios/feminicídios, 
29 lesões corporais seguidas de morte e 16 latrocínios = 416 registros - GEAC/DINE/SSP/SC - 2024.3*
"@
$ws.Range("E14").Value = $e14

$f14 = @"
SIM.

O Indicador e Fórmula de Cálculo informado foi:



4*Descrição do Indicador:
Taxa(%) de ocorrências atendidas.4*



5*Fórmula:
(Nº de ocorrências atendidas após a aquisição*100 / Nº de ocorrências atendidas realizadasantes da aquisção) - 100.5*

O indicador e a fórmula de cálculo são adequados para o eficiente monitoramento da meta.
"@
$ws.Range("F14").Value = $f14

$g14 = @"
SIM.

A Meta informada foi:



6*1ª Diretriz: Enfrentamento a criminalidade violenta. Ações Estratégicas 1. Redução dos homicídios. Taxa igual ou menor que 6,0/100 mil/hab. até 2030.6*



Existe aderência da referida Meta à Política Pública.
"@
$ws.Range("G14").Value = $g14

$h14 = @"
SIM.

A Meta informada foi:



7*Meta 1: Reduzir a taxa nacional de homicídios para abaixo de 16 mortes por 100 mil habitantes até 2030.7*



Existe aderência da referida Meta Específica à Política informada.
"@
$ws.Range("H14").Value = $h14

$i14 = @"
SIM.

A política informada foi:



8*Política de Enfrentamento da Criminalidade Violenta.8*



Existe aderência da referida Meta Específica à Política informada.
"@
$ws.Range("I14").Value = $i14

# --- Linha 25: era uma copia da linha 14 (placeholders 2*/3*/.../8*) e foi
# esvaziada, pois o conteudo real ja consta na linha 14 mesclada A14:A24 etc.
$ws.Range("A25:I25").Value = ""
